$d = $word.ActiveDocument
$word.Selection.Find.Execute("Java", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$word.Selection.TypeText("C#")

$rngC = $d.Range(244, 246)
$rngC.Select()
$word.Selection.LanguageID = 1036

$rngCode = $d.Range(236, 240)
Write-Output ("rngCode text: " + $rngCode.Text)
$rngCode.Select()
$word.Selection.LanguageID = 1036
